$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 2-3: mfd_hab1 (N) becomes "Sandfilter"; mfd_hab2 (O) and mfd_hab3 (P) are cleared (removed).
foreach ($r in 2..3) {
    $ws.Cells.Item($r, 14).Value = "Sandfilter"   # N: mfd_hab1
    $ws.Cells.Item($r, 15).ClearContents()        # O: mfd_hab2
    $ws.Cells.Item($r, 16).ClearContents()        # P: mfd_hab3
}

# Rows 4-7: habitat_typenumber (F) 6310 -> 3100; mfd_hab1 (N) Urban -> Other;
# mfd_hab2 (O) takes the old mfd_hab3 value "Landfill"; mfd_hab3 (P) is cleared.
foreach ($r in 4..7) {
    $fcell = $ws.Cells.Item($r, 6)
    $fcell.NumberFormat = "@"                     # keep habitat_typenumber as text
    $fcell.Value = "3100"                         # F: habitat_typenumber
    $ws.Cells.Item($r, 14).Value = "Other"        # N: mfd_hab1
    $ws.Cells.Item($r, 15).Value = "Landfill"     # O: mfd_hab2
    $ws.Cells.Item($r, 16).ClearContents()        # P: mfd_hab3
}

# Rows 8-29: habitat_typenumber (F) 2110 -> 2100; mfd_hab1 (N) Urban -> Other;
# mfd_hab2 (O) takes the old mfd_hab3 value "Landfill"; mfd_hab3 (P) is cleared.
foreach ($r in 8..29) {
    $fcell = $ws.Cells.Item($r, 6)
    $fcell.NumberFormat = "@"                     # keep habitat_typenumber as text
    $fcell.Value = "2100"                         # F: habitat_typenumber
    $ws.Cells.Item($r, 14).Value = "Other"        # N: mfd_hab1
    $ws.Cells.Item($r, 15).Value = "Landfill"     # O: mfd_hab2
    $ws.Cells.Item($r, 16).ClearContents()        # P: mfd_hab3
}
